$wb = $excel.ActiveWorkbook

function Set-RowData {
    param($ws, $r, $a, $b, $c, $d, $e, $f, $g, $h, $i)
    $ws.Cells.Item($r, 1).NumberFormat = "General"
    $ws.Cells.Item($r, 1).Value = $a
    $ws.Cells.Item($r, 2).NumberFormat = "@"
    $ws.Cells.Item($r, 2).Value = $b
    $ws.Cells.Item($r, 3).NumberFormat = "@"
    $ws.Cells.Item($r, 3).Value = $c
    $ws.Cells.Item($r, 4).NumberFormat = "@"
    $ws.Cells.Item($r, 4).Value = $d
    $ws.Cells.Item($r, 5).NumberFormat = "@"
    $ws.Cells.Item($r, 5).Value = $e
    $ws.Cells.Item($r, 6).NumberFormat = "General"
    $ws.Cells.Item($r, 6).Value = $f
    $ws.Cells.Item($r, 7).NumberFormat = "General"
    $ws.Cells.Item($r, 7).Value = $g
    $ws.Cells.Item($r, 8).NumberFormat = "@"
    $ws.Cells.Item($r, 8).Value = $h
    $ws.Cells.Item($r, 9).NumberFormat = "@"
    $ws.Cells.Item($r, 9).Value = $i
}

$ws = $wb.Worksheets.Item("展览")
Set-RowData $ws 2 "1" "2024-03-16" "合肥·CW国潮动漫游戏嘉年华" "南京路与庐州大道交汇处 合肥滨湖国际会展中心" "2024.03.16 09:30-03.17 17:00" 3024 65 "https://show.bilibili.com/platform/detail.html?id=81284" "//i0.hdslb.com/bfs/openplatform/202401/38B92fWF1705995243803.jpeg"
Set-RowData $ws 3 "2" "2024-03-17" "合肥·CW国潮动漫游戏嘉年华-赵路内场" "南京路与庐州大道交汇处 合肥滨湖国际会展中心" "2024.03.17 09:00-03.17 17:00" 732 "已售罄" "https://show.bilibili.com/platform/detail.html?id=81954" "//i0.hdslb.com/bfs/openplatform/202402/2PVn1ahm1708481741272.jpeg"
Set-RowData $ws 4 "3" "2024-03-23" "合肥·原&铁&崩 only展" "金寨路与天堂窄路交叉口 梵木艺术中心" "2024.03.23 09:00-03.23 17:00" 105 58 "https://show.bilibili.com/platform/detail.html?id=81574" "//i2.hdslb.com/bfs/openplatform/202401/0V5uyX6C1706697212904.png"
Set-RowData $ws 5 "4" "2024-04-04" "合肥· 第二届漫画城市动漫展 -故事再次开始" "凤淮路与固镇路西北角 庐阳全民健身中心" "2024.04.04 09:00-04.05 17:00" 6777 60 "https://show.bilibili.com/platform/detail.html?id=78898" "//i2.hdslb.com/bfs/openplatform/202402/3NgyB9761708333056023.jpeg"
Set-RowData $ws 6 "5" "2024-04-04" "合肥·环形宇宙动漫游戏嘉年华" "南京路与庐州大道交汇处 合肥滨湖国际会展中心" "2024.04.04 09:30-04.05 17:00" 1776 "不可售" "https://show.bilibili.com/platform/detail.html?id=81916" "//i0.hdslb.com/bfs/openplatform/202402/1lGzmBT61708336972816.jpeg"
Set-RowData $ws 7 "6" "2024-04-04" "合肥·环形宇宙动漫游戏嘉年华内场票-钱文青" "锦绣大道3899号 合肥滨湖会展中心" "2024.04.04 09:00-04.04 17:00" 7 238 "https://show.bilibili.com/platform/detail.html?id=82375" "//i1.hdslb.com/bfs/openplatform/202403/3cHtIycW1709692273438.jpeg"
Set-RowData $ws 8 "7" "2024-04-04" "合肥·第二届漫画城市动漫展内场-柯暮卿" "凤淮路与固镇路西北角 庐阳全民健身中心" "2024.04.04 10:00-04.04 17:00" 23 158 "https://show.bilibili.com/platform/detail.html?id=82192" "//i0.hdslb.com/bfs/openplatform/202402/tcAAj9aj1709193127773.jpeg"
Set-RowData $ws 9 "8" "2024-04-04" "合肥·第二届漫画城市动漫展内场-风袖" "凤淮路与固镇路西北角 庐阳全民健身中心" "2024.04.04 10:00-04.04 17:00" 31 158 "https://show.bilibili.com/platform/detail.html?id=82191" "//i2.hdslb.com/bfs/openplatform/202402/UZiEzBcc1709192469627.jpeg"
Set-RowData $ws 10 "9" "2024-04-05" "合肥· 第二届漫画城市动漫展内场-《琅声雅集》" "凤淮路与固镇路西北角 庐阳全民健身中心" "2024.04.05 10:00-04.05 17:00" 62 168 "https://show.bilibili.com/platform/detail.html?id=82189" "//i0.hdslb.com/bfs/openplatform/202402/ZKkx4hTN1709191842946.jpeg"
Set-RowData $ws 11 "10" "2024-04-21" "合肥·银魂only" "濉溪路118号 合肥栢景假日酒店" "2024.04.21 09:00-04.21 17:00" 127 55 "https://show.bilibili.com/platform/detail.html?id=82145" "//i1.hdslb.com/bfs/openplatform/202402/A0Tb5SQ51709091316985.jpeg"
Set-RowData $ws 12 "11" "2024-05-01" "合肥·Look Look动漫嘉年华" "新站区东方大道288号 少荃体育中心" "2024.05.01 10:00-05.01 17:30" 137 29.9 "https://show.bilibili.com/platform/detail.html?id=82311" "//i2.hdslb.com/bfs/openplatform/202403/jbUNtkAQ1709619599897.png"
Set-RowData $ws 13 "12" "2024-05-03" "合肥·BH动漫游戏展" "科技园路与葡萄园路交口包河区现代农业示范园8号 圩乐田园生态营地" "2024.05.03 10:00-05.04 16:00" 27 40 "https://show.bilibili.com/platform/detail.html?id=82199" "//i0.hdslb.com/bfs/openplatform/202402/cSR2xlY61709195356978.jpeg"
Set-RowData $ws 14 "13" "2024-05-18" "合肥·梦时空SPO1动漫展（取消）" "阜阳路16号 银瑞林国际大酒店" "2024.05.18 10:00-05.18 17:00" 132 "不可售" "https://show.bilibili.com/platform/detail.html?id=80207" "//i2.hdslb.com/bfs/openplatform/202312/tQQOHYE01703574162111.jpeg"

$ws = $wb.Worksheets.Item("全部类型")
Set-RowData $ws 2 "1" "2024-03-16" "合肥·CW国潮动漫游戏嘉年华" "南京路与庐州大道交汇处 合肥滨湖国际会展中心" "2024.03.16 09:30-03.17 17:00" 3024 65 "https://show.bilibili.com/platform/detail.html?id=81284" "//i0.hdslb.com/bfs/openplatform/202401/38B92fWF1705995243803.jpeg"
Set-RowData $ws 3 "2" "2024-03-17" "合肥·CW国潮动漫游戏嘉年华-赵路内场" "南京路与庐州大道交汇处 合肥滨湖国际会展中心" "2024.03.17 09:00-03.17 17:00" 732 "已售罄" "https://show.bilibili.com/platform/detail.html?id=81954" "//i0.hdslb.com/bfs/openplatform/202402/2PVn1ahm1708481741272.jpeg"
Set-RowData $ws 4 "3" "2024-03-23" "合肥·原&铁&崩 only展" "金寨路与天堂窄路交叉口 梵木艺术中心" "2024.03.23 09:00-03.23 17:00" 105 58 "https://show.bilibili.com/platform/detail.html?id=81574" "//i2.hdslb.com/bfs/openplatform/202401/0V5uyX6C1706697212904.png"
Set-RowData $ws 5 "4" "2024-04-04" "合肥· 第二届漫画城市动漫展 -故事再次开始" "凤淮路与固镇路西北角 庐阳全民健身中心" "2024.04.04 09:00-04.05 17:00" 6777 60 "https://show.bilibili.com/platform/detail.html?id=78898" "//i2.hdslb.com/bfs/openplatform/202402/3NgyB9761708333056023.jpeg"
Set-RowData $ws 6 "5" "2024-04-04" "合肥·环形宇宙动漫游戏嘉年华" "南京路与庐州大道交汇处 合肥滨湖国际会展中心" "2024.04.04 09:30-04.05 17:00" 1776 "不可售" "https://show.bilibili.com/platform/detail.html?id=81916" "//i0.hdslb.com/bfs/openplatform/202402/1lGzmBT61708336972816.jpeg"
Set-RowData $ws 7 "6" "2024-04-04" "合肥·环形宇宙动漫游戏嘉年华内场票-钱文青" "锦绣大道3899号 合肥滨湖会展中心" "2024.04.04 09:00-04.04 17:00" 7 238 "https://show.bilibili.com/platform/detail.html?id=82375" "//i1.hdslb.com/bfs/openplatform/202403/3cHtIycW1709692273438.jpeg"
Set-RowData $ws 8 "7" "2024-04-04" "合肥·第二届漫画城市动漫展内场-柯暮卿" "凤淮路与固镇路西北角 庐阳全民健身中心" "2024.04.04 10:00-04.04 17:00" 23 158 "https://show.bilibili.com/platform/detail.html?id=82192" "//i0.hdslb.com/bfs/openplatform/202402/tcAAj9aj1709193127773.jpeg"
Set-RowData $ws 9 "8" "2024-04-04" "合肥·第二届漫画城市动漫展内场-风袖" "凤淮路与固镇路西北角 庐阳全民健身中心" "2024.04.04 10:00-04.04 17:00" 31 158 "https://show.bilibili.com/platform/detail.html?id=82191" "//i2.hdslb.com/bfs/openplatform/202402/UZiEzBcc1709192469627.jpeg"
Set-RowData $ws 10 "9" "2024-04-05" "合肥· 第二届漫画城市动漫展内场-《琅声雅集》" "凤淮路与固镇路西北角 庐阳全民健身中心" "2024.04.05 10:00-04.05 17:00" 62 168 "https://show.bilibili.com/platform/detail.html?id=82189" "//i0.hdslb.com/bfs/openplatform/202402/ZKkx4hTN1709191842946.jpeg"
Set-RowData $ws 11 "10" "2024-04-21" "合肥·银魂only" "濉溪路118号 合肥栢景假日酒店" "2024.04.21 09:00-04.21 17:00" 127 55 "https://show.bilibili.com/platform/detail.html?id=82145" "//i1.hdslb.com/bfs/openplatform/202402/A0Tb5SQ51709091316985.jpeg"
Set-RowData $ws 12 "11" "2024-05-01" "合肥·Look Look动漫嘉年华" "新站区东方大道288号 少荃体育中心" "2024.05.01 10:00-05.01 17:30" 137 29.9 "https://show.bilibili.com/platform/detail.html?id=82311" "//i2.hdslb.com/bfs/openplatform/202403/jbUNtkAQ1709619599897.png"
Set-RowData $ws 13 "12" "2024-05-03" "合肥·BH动漫游戏展" "科技园路与葡萄园路交口包河区现代农业示范园8号 圩乐田园生态营地" "2024.05.03 10:00-05.04 16:00" 27 40 "https://show.bilibili.com/platform/detail.html?id=82199" "//i0.hdslb.com/bfs/openplatform/202402/cSR2xlY61709195356978.jpeg"
Set-RowData $ws 14 "13" "2024-05-18" "合肥·梦时空SPO1动漫展（取消）" "阜阳路16号 银瑞林国际大酒店" "2024.05.18 10:00-05.18 17:00" 132 "不可售" "https://show.bilibili.com/platform/detail.html?id=80207" "//i2.hdslb.com/bfs/openplatform/202312/tQQOHYE01703574162111.jpeg"

